# aggiunta indicatori italia modulo CR SISTEMA
# Adds 11 new "CRSYS_*" indicator rows to the "r AnalysisUnit_Variable" sheet,
# mirroring the layout/formatting of the existing rows (columns A/B/C/E/F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r AnalysisUnit_Variable")

$newNames = @(
    "CRSYS_TUA",
    "CRSYS_TUR",
    "CRSYS_TUS",
    "CRSYS_UCFBT",
    "CRSYS_Q_DER",
    "CRSYS_TUC",
    "CRSYS_ITUR",
    "CRSYS_UA",
    "CRSYS_US",
    "CRSYS_UR",
    "CRSYS_UT"
)

$startRow = 42
$lastExistingRow = 41

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $row = $startRow + $i
    $name = $newNames[$i]

    $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 5).Value = "CUSTOMER"
    $ws.Cells.Item($row, 6).Value = $name

    # Match formatting of the last pre-existing data row (row 41: B/C/F styled).
    $ws.Cells.Item($row, 2).Style = $ws.Cells.Item($lastExistingRow, 2).Style
    $ws.Cells.Item($row, 3).Style = $ws.Cells.Item($lastExistingRow, 3).Style
    $ws.Cells.Item($row, 6).Style = $ws.Cells.Item($lastExistingRow, 6).Style
}

$endRow = $startRow + $newNames.Count - 1

# Update the view to reflect the newly added rows, as in the target workbook
# (matches the saved selection/scroll position captured in the authored diff).
$ws.Activate()
$ws.Range("F47").Select()
